# Applies the edit described by the commit "Fruta / hortaliza, semanal":
# A new day's price record is inserted at the top of the data (row 24),
# pushing the existing rows 24-132 down to 25-133, and the worksheet's
# dimension grows from A1:R132 to A1:R133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right before the current row 24; this shifts rows
# 24..132 down to 25..133 (carrying along all of their existing values
# and formatting), and extends the sheet dimension automatically.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new daily record.
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44676
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112001
$ws.Cells.Item(24, 7).Value = "Berenjena"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 2600
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 9500
$ws.Cells.Item(24, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(24, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24, 16).Value = 190
$ws.Cells.Item(24, 17).Value = 50
$ws.Cells.Item(24, 18).Value = "Hortaliza"
